$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The report header block (previously in columns C/D, rows 2-4) moves one
# column to the left (B/C) and one row down (rows 3-4), and a brand new
# "Ruta:" / "DEFAULT" row is inserted above it at row 2. The old stray
# formatted-but-empty cell at C4 disappears entirely.

# Wipe out the old header block (contents + formatting) so nothing is left
# behind in column D or in the vacated cells.
$ws.Range("C2:D4").Clear()

# New row 2: Ruta / DEFAULT
$ws.Range("B2").Value = "Ruta:"
$ws.Range("C2").Value = "DEFAULT"

# Row 3: Fecha (date) - shifted from the old C2/D2
$ws.Range("B3").Value = "Fecha : "
$ws.Range("C3").Value = 43692
$ws.Range("C3").NumberFormat = "d-mmm-yy"

# Row 4: Repartidor / Nombre de repartidor - shifted from the old C3/D3
$ws.Range("B4").Value = "Repartidor: "
$ws.Range("C4").Value = "Nombre de repartidor"

# Move the active selection like the edit author left it
$ws.Range("D6").Select() | Out-Null
